$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2025-09-30"
$ws.Range("A49").Style = "Normal"

$ws.Range("B49").Value = "15:19:33"
$ws.Range("C49").Value = "1.00 EUR = 1,638.2586"
